# Fix The Burning Kor reminder cards
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also updates the ExternalData_1 defined name reference automatically)
$ws.Name = "Nemesis Cards"

# Clear A4 (was a duplicate/empty-string value, now genuinely blank)
$ws.Range("A4").ClearContents()

# Update the two "Burning Kor" reminder-card cells with corrected/expanded rule text
$ws.Range("K669").Value = '>>>- Once per turn, during any player''s casting or main phase, that player may move their mage token to an adjacent space for free.#- During any player''s casting or main phase, that player may discard a card to move their mage token to an adjacent space.#- Mage tokens cannot occupy the same space as embers or other mage tokens.#- Nothing can move onto Gravehold.#- Players can''t move their mage tokens into any of the ignition points.'
$ws.Range("K670").Value = '>>>- Players can only deal damage to embers that are orthogonally adjacent to their mage tokens.#- When a player deals damage to an ember, they can distribute that damage among any number of embers in spaces adjacent to their mage token.#- Embers are minions and have 1 life.#- Reduce to 0 all damage dealt to The Burning Kor by players not on the edge of the board.'

# Restore selection/scroll state to match the saved view
$ws.Range("K671").Select()
